$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("export")
$ws.Range("G1").Interior.Color = 65535
$ws.Range("G2").Interior.Color = 5296274
$ws.Range("G3").Interior.ThemeColor = 2
$ws.Range("G4").Interior.Color = 255
Write-Output "done"
